$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14369  # was 14362
$ws.Range("F4").Value = 14519  # was 14500
$ws.Range("F6").Value = 1413  # was 1412
$ws.Range("F7").Value = 5936  # was 5934
$ws.Range("F12").Value = 201  # was 200
$ws.Range("F13").Value = 1566  # was 1565
$ws.Range("F15").Value = 2135  # was 2134
$ws.Range("F16").Value = 1223  # was 1221
$ws.Range("F17").Value = 1870  # was 1868
$ws.Range("F22").Value = 829  # was 828
$ws.Range("F23").Value = 3389  # was 3384
$ws.Range("F25").Value = 319  # was 317
$ws.Range("F26").Value = 2448  # was 2445
$ws.Range("F30").Value = 1819  # was 1816
$ws.Range("F32").Value = 1439  # was 1431
$ws.Range("F35").Value = 4968  # was 4962
$ws.Range("F36").Value = 4933  # was 4928
$ws.Range("F41").Value = 3315  # was 3313
$ws.Range("F45").Value = 118  # was 117
$ws.Range("F46").Value = 98  # was 96
$ws.Range("F48").Value = 625  # was 621
$ws.Range("F49").Value = 302  # was 301

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 127  # was 126
$ws.Range("F26").Value = 70  # was 69

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7677  # was 7669
$ws.Range("F3").Value = 255  # was 253
$ws.Range("F4").Value = 875  # was 869

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7677  # was 7669
$ws.Range("F4").Value = 255  # was 253
$ws.Range("F5").Value = 875  # was 869
$ws.Range("F7").Value = 14369  # was 14362
$ws.Range("F8").Value = 14519  # was 14500
$ws.Range("F10").Value = 1413  # was 1412
$ws.Range("F11").Value = 5936  # was 5934
$ws.Range("F13").Value = 127  # was 126
$ws.Range("F16").Value = 1566  # was 1565
$ws.Range("F19").Value = 829  # was 828
$ws.Range("F20").Value = 3389  # was 3384
$ws.Range("F21").Value = 319  # was 317
$ws.Range("F22").Value = 2448  # was 2445
$ws.Range("F25").Value = 1819  # was 1816
$ws.Range("F32").Value = 1439  # was 1431
$ws.Range("F35").Value = 4968  # was 4962
$ws.Range("F36").Value = 4933  # was 4928
$ws.Range("F39").Value = 3315  # was 3313
$ws.Range("F42").Value = 118  # was 117
$ws.Range("F44").Value = 98  # was 96
$ws.Range("F45").Value = 625  # was 621
$ws.Range("F46").Value = 302  # was 301
$ws.Range("F47").Value = 70  # was 69
